$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update values per the diff
$ws.Range("B1").Value = 4
$ws.Range("B2").Value = 10
$ws.Range("B6").Value = 2
$ws.Range("B7").Value = 4

# Update the active cell selection to D15
$ws.Range("D15").Select()
